# Apply updated transition-probability values to the team-specific matrix
# (Northwest (WA)_B) worksheet. The sheet is a row-normalized state
# transition matrix; the cells below were populated with newly simulated
# probabilities (previously 0 placeholders) as part of adding more games,
# speeding up the simulate-game logic, and drafting optimization logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2"  = 0.5
    "S2"  = 0.5

    "S3"  = 1

    "F6"  = 0.1111111111111111
    "J6"  = 0.2222222222222222
    "Q6"  = 0.1111111111111111
    "R6"  = 0.2222222222222222
    "S6"  = 0.3333333333333333

    "O7"  = 0.4
    "S7"  = 0.6

    "B8"  = 0.03333333333333333
    "F8"  = 0.03333333333333333
    "J8"  = 0.03333333333333333
    "O8"  = 0.03333333333333333
    "Q8"  = 0.1333333333333333
    "R8"  = 0.1
    "S8"  = 0.6333333333333333

    "O9"  = 0.2
    "Q9"  = 0.1
    "R9"  = 0.1
    "S9"  = 0.6

    "B10" = 0.03448275862068965
    "F10" = 0.1379310344827586
    "O10" = 0.06896551724137931
    "Q10" = 0.2068965517241379
    "R10" = 0.03448275862068965
    "S10" = 0.5172413793103449

    "G11" = 0.2727272727272727
    "J11" = 0.09090909090909091
    "K11" = 0.3636363636363636
    "L11" = 0.2727272727272727

    "G12" = 0.6666666666666666
    "J12" = 0.3333333333333333

    "H15" = 0.2941176470588235
    "I15" = 0.05882352941176471
    "J15" = 0.1176470588235294
    "K15" = 0.1176470588235294
    "O15" = 0.1176470588235294
    "S15" = 0.2941176470588235

    "H17" = 0.1818181818181818
    "I17" = 0.1818181818181818
    "J17" = 0.2727272727272727
    "K17" = 0.09090909090909091
    "O17" = 0.1818181818181818
    "S17" = 0.09090909090909091

    "F18" = 0.1428571428571428
    "I18" = 0.1428571428571428
    "J18" = 0.5714285714285714
    "O18" = 0.1428571428571428

    "F19" = 0.03278688524590164
    "H19" = 0.3770491803278688
    "I19" = 0.09836065573770492
    "J19" = 0.2459016393442623
    "K19" = 0.06557377049180328
    "O19" = 0.06557377049180328
    "S19" = 0.1147540983606557
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
